# The workbook moved from the "Boils_per_Charge" folder into "Components"
# as part of merging updated Excel files into folders. Reset the input
# parameters (Battery_Capacity, Watts, T_Boil) to placeholder values of 1;
# the Boils_per_Charge formula (B6 = B3/(B4*B5)) recalculates automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("Battery_Capacity").Value = 1
$ws.Range("Watts").Value = 1
$ws.Range("T_Boil").Value = 1
